# Updates cryptos list data (prices and hourly volume percentages) to reflect
# the latest GitHub Actions scrape. Also accounts for a reordering of two
# rows (Monero / LidoDAOToken swapped positions with updated data).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '45.349.64'
$ws.Range("E2").Value = '  +4.11%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.425.67'
$ws.Range("E3").Value = '  +0.11%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '318.29'
$ws.Range("E5").Value = '  +3.81%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '102.52'
$ws.Range("E6").Value = '  +5.27%  '
$ws.Range("E7").Value = '  +1.45%  '
$ws.Range("E8").Value = '  -0.07%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.531'
$ws.Range("E9").Value = '  +8.02%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.63'
$ws.Range("E10").Value = '  +1.36%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0804'
$ws.Range("E11").Value = '  +0.73%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.13'
$ws.Range("E13").Value = '  -2.04%  '
$ws.Range("E14").Value = '  +2.05%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.804.11'
$ws.Range("E15").Value = '  +0.59%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.437.26'
$ws.Range("E16").Value = '  +1.00%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.844'
$ws.Range("E17").Value = '  +2.08%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '45.231.72'
$ws.Range("E18").Value = '  +3.89%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.24'
$ws.Range("E19").Value = '  +1.44%  '
$ws.Range("E20").Value = '  -0.99%  '
$ws.Range("E21").Value = '  +2.25%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '68.82'
$ws.Range("E22").Value = '  +1.05%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '244.23'
$ws.Range("E23").Value = '  +2.72%  '
$ws.Range("E24").Value = '  -0.02%  '
$ws.Range("E25").Value = '  +1.84%  '
$ws.Range("E26").Value = '  -0.03%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '25.58'
$ws.Range("E27").Value = '  +2.31%  '
$ws.Range("E28").Value = '  -6.78%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.59'
$ws.Range("E29").Value = '  +1.41%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '49.09'
$ws.Range("E30").Value = '  +2.02%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '32.91'
$ws.Range("E31").Value = '  +1.89%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.127'
$ws.Range("E32").Value = '  +5.08%  '
$ws.Range("E33").Value = '  +9.44%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.21'
$ws.Range("E34").Value = '  +1.50%  '
$ws.Range("E35").Value = '  +0.31%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0764'
$ws.Range("E36").Value = '  +1.55%  '
$ws.Range("E37").Value = '  -1.30%  '
$ws.Range("E38").Value = '  +1.64%  '
$ws.Range("B39").Value = 'LidoDAOToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.86'
$ws.Range("E39").Value = '  -2.14%  '
$ws.Range("B40").Value = 'Monero'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '126.04'
$ws.Range("E40").Value = '  -4.32%  '
$ws.Range("E41").Value = '  -2.45%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.110'
$ws.Range("E42").Value = '  +1.01%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '20.63'
$ws.Range("E43").Value = '  -1.93%  '
$ws.Range("E44").Value = '  +2.52%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.935.22'
$ws.Range("E45").Value = '  -0.51%  '
$ws.Range("E46").Value = '  -2.95%  '
$ws.Range("E47").Value = '  +3.14%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.83'
$ws.Range("E48").Value = '  +16.16%  '
$ws.Range("E49").Value = '  -2.33%  '
$ws.Range("E50").Value = '  +5.54%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '53.91'
$ws.Range("E51").Value = '  +2.26%  '
